$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("E2").Value = 3
$ws.Range("G2").Value = 21.443826
$ws.Range("H2").Value = 64.331478
$ws.Range("I2").Value = 0.6062978927103765
$ws.Range("J2").Value = 0.6062978927103765
$ws.Range("K2").Value = 3
$ws.Range("M2").Value = 10.90224333333333
$ws.Range("N2").Value = 32.70673
$ws.Range("O2").Value = 0.1536506229779223
$ws.Range("P2").Value = 0.1536506229779223
$ws.Range("Q2").Value = 233.78580904966
$ws.Range("R2").Value = 2104.07228144694
$ws.Range("S2").Value = 0.09315804892515084
$ws.Range("T2").Value = 0.09315804892515085

$ws.Range("E3").Value = 3
$ws.Range("G3").Value = 21.443826
$ws.Range("H3").Value = 64.331478
$ws.Range("I3").Value = 0.6062978927103765
$ws.Range("J3").Value = 0.6062978927103765
$ws.Range("K3").Value = 3
$ws.Range("M3").Value = 27.20435933333333
$ws.Range("N3").Value = 81.613078
$ws.Range("O3").Value = 0.3834042803375869
$ws.Range("P3").Value = 0.3834042803375869
$ws.Range("Q3").Value = 583.365547985476
$ws.Range("R3").Value = 5250.289931869284
$ws.Range("S3").Value = 0.2324572072248174
$ws.Range("T3").Value = 0.2324572072248174

$ws.Range("E4").Value = 3
$ws.Range("G4").Value = 21.443826
$ws.Range("H4").Value = 64.331478
$ws.Range("I4").Value = 0.6062978927103765
$ws.Range("J4").Value = 0.6062978927103765
$ws.Range("K4").Value = 3
$ws.Range("M4").Value = 5.916218
$ws.Range("N4").Value = 17.748654
$ws.Range("O4").Value = 0.08338014054353927
$ws.Range("P4").Value = 0.08338014054353927
$ws.Range("Q4").Value = 126.866349370068
$ws.Range("R4").Value = 1141.797144330612
$ws.Range("S4").Value = 0.05055320350544289
$ws.Range("T4").Value = 0.05055320350544289

$ws.Range("E5").Value = 3
$ws.Range("G5").Value = 21.443826
$ws.Range("H5").Value = 64.331478
$ws.Range("I5").Value = 0.6062978927103765
$ws.Range("J5").Value = 0.6062978927103765
$ws.Range("K5").Value = 3
$ws.Range("M5").Value = 6.012963333333334
$ws.Range("N5").Value = 18.03889
$ws.Range("O5").Value = 0.08474361962599786
$ws.Range("P5").Value = 0.08474361962599786
$ws.Range("Q5").Value = 128.94093946438
$ws.Range("R5").Value = 1160.46845517942
$ws.Range("S5").Value = 0.0513798779998922
$ws.Range("T5").Value = 0.0513798779998922

$ws.Range("E6").Value = 3
$ws.Range("G6").Value = 21.443826
$ws.Range("H6").Value = 64.331478
$ws.Range("I6").Value = 0.6062978927103765
$ws.Range("J6").Value = 0.6062978927103765
$ws.Range("K6").Value = 3
$ws.Range("M6").Value = 3.937401666666666
$ws.Range("N6").Value = 11.812205
$ws.Range("O6").Value = 0.05549171858491902
$ws.Range("P6").Value = 0.05549171858491902
$ws.Range("Q6").Value = 84.43295623211
$ws.Range("R6").Value = 759.89660608899
$ws.Range("S6").Value = 0.03364451204091364
$ws.Range("T6").Value = 0.03364451204091364

$ws.Range("E7").Value = 3
$ws.Range("G7").Value = 21.443826
$ws.Range("H7").Value = 64.331478
$ws.Range("I7").Value = 0.6062978927103765
$ws.Range("J7").Value = 0.6062978927103765
$ws.Range("K7").Value = 3
$ws.Range("M7").Value = 16.981576
$ws.Range("N7").Value = 50.944728
$ws.Range("O7").Value = 0.2393296179300346
$ws.Range("P7").Value = 0.2393296179300347
$ws.Range("Q7").Value = 364.149960949776
$ws.Range("R7").Value = 3277.349648547985
$ws.Range("S7").Value = 0.1451050430141595
$ws.Range("T7").Value = 0.1451050430141596

$ws.Range("E8").Value = 3
$ws.Range("G8").Value = 3.641794
$ws.Range("H8").Value = 10.925382
$ws.Range("I8").Value = 0.1029672609675761
$ws.Range("J8").Value = 0.1029672609675761
$ws.Range("K8").Value = 3
$ws.Range("M8").Value = 10.90224333333333
$ws.Range("N8").Value = 32.70673
$ws.Range("O8").Value = 0.1536506229779223
$ws.Range("P8").Value = 0.1536506229779223
$ws.Range("Q8").Value = 39.70372435787333
$ws.Range("R8").Value = 357.33351922086
$ws.Range("S8").Value = 0.01582098379399836
$ws.Range("T8").Value = 0.01582098379399837

$ws.Range("E9").Value = 3
$ws.Range("G9").Value = 3.641794
$ws.Range("H9").Value = 10.925382
$ws.Range("I9").Value = 0.1029672609675761
$ws.Range("J9").Value = 0.1029672609675761
$ws.Range("K9").Value = 3
$ws.Range("M9").Value = 27.20435933333333
$ws.Range("N9").Value = 81.613078
$ws.Range("O9").Value = 0.3834042803375869
$ws.Range("P9").Value = 0.3834042803375869
$ws.Range("Q9").Value = 99.07267259397734
$ws.Range("R9").Value = 891.6540533457961
$ws.Range("S9").Value = 0.039478088589606
$ws.Range("T9").Value = 0.03947808858960601

$ws.Range("E10").Value = 3
$ws.Range("G10").Value = 3.641794
$ws.Range("H10").Value = 10.925382
$ws.Range("I10").Value = 0.1029672609675761
$ws.Range("J10").Value = 0.1029672609675761
$ws.Range("K10").Value = 3
$ws.Range("M10").Value = 5.916218
$ws.Range("N10").Value = 17.748654
$ws.Range("O10").Value = 0.08338014054353927
$ws.Range("P10").Value = 0.08338014054353927
$ws.Range("Q10").Value = 21.545647215092
$ws.Range("R10").Value = 193.910824935828
$ws.Range("S10").Value = 0.008585424690859779
$ws.Range("T10").Value = 0.008585424690859779

$ws.Range("E11").Value = 3
$ws.Range("G11").Value = 3.641794
$ws.Range("H11").Value = 10.925382
$ws.Range("I11").Value = 0.1029672609675761
$ws.Range("J11").Value = 0.1029672609675761
$ws.Range("K11").Value = 3
$ws.Range("M11").Value = 6.012963333333334
$ws.Range("N11").Value = 18.03889
$ws.Range("O11").Value = 0.08474361962599786
$ws.Range("P11").Value = 0.08474361962599786
$ws.Range("Q11").Value = 21.89797378955334
$ws.Range("R11").Value = 197.08176410598
$ws.Range("S11").Value = 0.008725818397367122
$ws.Range("T11").Value = 0.008725818397367122

$ws.Range("E12").Value = 3
$ws.Range("G12").Value = 3.641794
$ws.Range("H12").Value = 10.925382
$ws.Range("I12").Value = 0.1029672609675761
$ws.Range("J12").Value = 0.1029672609675761
$ws.Range("K12").Value = 3
$ws.Range("M12").Value = 3.937401666666666
$ws.Range("N12").Value = 11.812205
$ws.Range("O12").Value = 0.05549171858491902
$ws.Range("P12").Value = 0.05549171858491902
$ws.Range("Q12").Value = 14.33920576525667
$ws.Range("R12").Value = 129.05285188731
$ws.Range("S12").Value = 0.005713830269072647
$ws.Range("T12").Value = 0.005713830269072647

$ws.Range("E13").Value = 3
$ws.Range("G13").Value = 3.641794
$ws.Range("H13").Value = 10.925382
$ws.Range("I13").Value = 0.1029672609675761
$ws.Range("J13").Value = 0.1029672609675761
$ws.Range("K13").Value = 3
$ws.Range("M13").Value = 16.981576
$ws.Range("N13").Value = 50.944728
$ws.Range("O13").Value = 0.2393296179300346
$ws.Range("P13").Value = 0.2393296179300347
$ws.Range("Q13").Value = 61.84340158734401
$ws.Range("R13").Value = 556.5906142860961
$ws.Range("S13").Value = 0.02464311522667215
$ws.Range("T13").Value = 0.02464311522667215

$ws.Range("E14").Value = 3
$ws.Range("G14").Value = 10.28284533333333
$ws.Range("H14").Value = 30.848536
$ws.Range("I14").Value = 0.2907348463220475
$ws.Range("J14").Value = 0.2907348463220475
$ws.Range("K14").Value = 3
$ws.Range("M14").Value = 10.90224333333333
$ws.Range("N14").Value = 32.70673
$ws.Range("O14").Value = 0.1536506229779223
$ws.Range("P14").Value = 0.1536506229779223
$ws.Range("Q14").Value = 112.1060819830311
$ws.Range("R14").Value = 1008.95473784728
$ws.Range("S14").Value = 0.0446715902587731
$ws.Range("T14").Value = 0.04467159025877311

$ws.Range("E15").Value = 3
$ws.Range("G15").Value = 10.28284533333333
$ws.Range("H15").Value = 30.848536
$ws.Range("I15").Value = 0.2907348463220475
$ws.Range("J15").Value = 0.2907348463220475
$ws.Range("K15").Value = 3
$ws.Range("M15").Value = 27.20435933333333
$ws.Range("N15").Value = 81.613078
$ws.Range("O15").Value = 0.3834042803375869
$ws.Range("P15").Value = 0.3834042803375869
$ws.Range("Q15").Value = 279.7382194170897
$ws.Range("R15").Value = 2517.643974753808
$ws.Range("S15").Value = 0.1114689845231636
$ws.Range("T15").Value = 0.1114689845231636

$ws.Range("E16").Value = 3
$ws.Range("G16").Value = 10.28284533333333
$ws.Range("H16").Value = 30.848536
$ws.Range("I16").Value = 0.2907348463220475
$ws.Range("J16").Value = 0.2907348463220475
$ws.Range("K16").Value = 3
$ws.Range("M16").Value = 5.916218
$ws.Range("N16").Value = 17.748654
$ws.Range("O16").Value = 0.08338014054353927
$ws.Range("P16").Value = 0.08338014054353927
$ws.Range("Q16").Value = 60.83555465228266
$ws.Range("R16").Value = 547.5199918705439
$ws.Range("S16").Value = 0.02424151234723662
$ws.Range("T16").Value = 0.02424151234723662

$ws.Range("E17").Value = 3
$ws.Range("G17").Value = 10.28284533333333
$ws.Range("H17").Value = 30.848536
$ws.Range("I17").Value = 0.2907348463220475
$ws.Range("J17").Value = 0.2907348463220475
$ws.Range("K17").Value = 3
$ws.Range("M17").Value = 6.012963333333334
$ws.Range("N17").Value = 18.03889
$ws.Range("O17").Value = 0.08474361962599786
$ws.Range("P17").Value = 0.08474361962599786
$ws.Range("Q17").Value = 61.83037195167112
$ws.Range("R17").Value = 556.4733475650401
$ws.Range("S17").Value = 0.02463792322873854
$ws.Range("T17").Value = 0.02463792322873854

$ws.Range("E18").Value = 3
$ws.Range("G18").Value = 10.28284533333333
$ws.Range("H18").Value = 30.848536
$ws.Range("I18").Value = 0.2907348463220475
$ws.Range("J18").Value = 0.2907348463220475
$ws.Range("K18").Value = 3
$ws.Range("M18").Value = 3.937401666666666
$ws.Range("N18").Value = 11.812205
$ws.Range("O18").Value = 0.05549171858491902
$ws.Range("P18").Value = 0.05549171858491902
$ws.Range("Q18").Value = 40.48769235354222
$ws.Range("R18").Value = 364.3892311818799
$ws.Range("S18").Value = 0.01613337627493274
$ws.Range("T18").Value = 0.01613337627493274

$ws.Range("E19").Value = 3
$ws.Range("G19").Value = 10.28284533333333
$ws.Range("H19").Value = 30.848536
$ws.Range("I19").Value = 0.2907348463220475
$ws.Range("J19").Value = 0.2907348463220475
$ws.Range("K19").Value = 3
$ws.Range("M19").Value = 16.981576
$ws.Range("N19").Value = 50.944728
$ws.Range("O19").Value = 0.2393296179300346
$ws.Range("P19").Value = 0.2393296179300347
$ws.Range("Q19").Value = 174.6189195242453
$ws.Range("R19").Value = 1571.570275718208
$ws.Range("S19").Value = 0.06958145968920297
$ws.Range("T19").Value = 0.06958145968920298

